$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate raw-data table (rows 16-27) that is no longer needed
$ws.Range("A16:G27").EntireRow.Delete()

# Update the Archaea/Other row (E12): orgs_with_bchlD_fs count increases by one
# because the MoxR family ATPase query was added, contributing one more match
$ws.Range("E12").Value = 47

# Recalculate so the SUM formula in row 13 reflects the new total
$excel.CalculateFull()

# Update the view state to match what was left on screen after the edit
$ws.Activate()
$ws.Range("B15").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

Write-Host ("Dimension check - UsedRange: " + $ws.UsedRange.Address())
Write-Host ("E12: " + $ws.Range("E12").Value())
Write-Host ("E13: " + $ws.Range("E13").Value())
